$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Row 2 values
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 90
$ws.Range("D2").Value = 25
$ws.Range("E2").Value = 30
$ws.Range("F2").Value = 120
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 36
$ws.Range("I2").Value = 82
$ws.Range("J2").Value = 28
$ws.Range("K2").Value = 44

# Update Row 3 values
$ws.Range("G3").Value = 170
$ws.Range("J3").Value = 146
$ws.Range("K3").Value = 389
$ws.Range("L3").Value = 19

# Update Row 4 values
$ws.Range("C4").Value = 235
$ws.Range("E4").Value = 140
$ws.Range("I4").Value = 167
$ws.Range("J4").Value = 20
$ws.Range("L4").Value = 30

# Update selection to match the final state (active cell K4, single-cell selection)
$ws.Range("K4").Select()
